$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency market data as of 2024-04-14
# Each entry: Row, Ticker (B), Name (C), Price (D), Market Cap (E), Volume (F), Change 24h (G)
$data = @(
    @(2, 'BTC', 'Bitcoin', 63788, 1251276371121, 64058889462, -5.5655),
    @(3, 'ETH', 'Ethereum', 3042.66, 363088105363, 34549843049, -6.75679),
    @(4, 'USDT', 'Tether', 0.999846, 107281972564, 125889492954, -0.01726),
    @(5, 'BNB', 'BNB', 548.6, 83824126043, 4037942392, -7.08988),
    @(6, 'SOL', 'Solana', 138.28, 61847455906, 12080316552, -8.052099999999999),
    @(7, 'USDC', 'USDC', 1, 32358687395, 11684867390, -0.0554),
    @(8, 'STETH', 'Lido Staked Ether', 3041.54, 28570464641, 408038060, -6.50143),
    @(9, 'XRP', 'XRP', 0.485721, 26752711942, 4697776472, -10.65718),
    @(10, 'TON', 'Toncoin', 6.44, 22321043075, 1026831209, -3.85049),
    @(11, 'DOGE', 'Dogecoin', 0.152334, 21675964275, 6110129035, -10.97325),
    @(12, 'ADA', 'Cardano', 0.456654, 15964033014, 1474309600, -9.4658),
    @(13, 'AVAX', 'Avalanche', 35.7, 13330755922, 1971060110, -6.73903),
    @(14, 'SHIB', 'Shiba Inu', 0.00002173, 12816697217, 1911392711, -12.03616),
    @(15, 'WBTC', 'Wrapped Bitcoin', 63924, 9930759518, 764761889, -5.47332),
    @(16, 'TRX', 'TRON', 0.110584, 9707030759, 794482894, -3.09989),
    @(17, 'BCH', 'Bitcoin Cash', 482.88, 9353845550, 1325774502, -8.88552),
    @(18, 'DOT', 'Polkadot', 6.55, 8876488618, 859864388, -8.03687),
    @(19, 'LINK', 'Chainlink', 13.61, 7917179658, 1402412091, -9.0284),
    @(20, 'MATIC', 'Polygon', 0.681348, 6298508280, 1188473538, -9.85017),
    @(21, 'LTC', 'Litecoin', 78.09999999999999, 5820436033, 1499889837, -8.69842),
    @(22, 'ICP', 'Internet Computer', 12.44, 5730989972, 418637263, -7.74713),
    @(23, 'NEAR', 'NEAR Protocol', 5.37, 5687923057, 1763259297, -6.06506),
    @(24, 'LEO', 'LEO Token', 5.92, 5474324308, 2165961, 2.61709),
    @(25, 'UNI', 'Uniswap', 7.11, 5342030404, 631267616, -9.83084),
    @(26, 'DAI', 'Dai', 1, 5119768598, 798888851, 0.07745),
    @(27, 'FDUSD', 'First Digital USD', 1.001, 3941873828, 11949300894, 0.00663),
    @(28, 'APT', 'Aptos', 9.25, 3873352095, 513013769, -8.01047),
    @(29, 'TAO', 'Bittensor', 572.95, 3805666688, 129472649, 11.94782),
    @(30, 'STX', 'Stacks', 2.63, 3775178693, 283197803, -1.67814),
    @(31, 'ETC', 'Ethereum Classic', 25.65, 3736131644, 651256629, -11.88141),
    @(32, 'MNT', 'Mantle', 1.1, 3610229282, 148717178, -7.60669),
    @(33, 'CRO', 'Cronos', 0.128259, 3396305543, 49025328, -6.68048),
    @(34, 'ATOM', 'Cosmos Hub', 8.19, 3179428255, 889703588, -8.401870000000001),
    @(35, 'XLM', 'Stellar', 0.108128, 3116808681, 408562307, -4.55396),
    @(36, 'FIL', 'Filecoin', 5.84, 3115527267, 789450224, -11.70413),
    @(37, 'OKB', 'OKB', 51.73, 3102861366, 27662462, -2.70965),
    @(38, 'VET', 'VeChain', 0.04098192, 2985566543, 196757479, -6.81483),
    @(39, 'RNDR', 'Render', 7.69, 2936613852, 550481339, -4.16378),
    @(40, 'ARB', 'Arbitrum', 1.12, 2932869705, 1187221157, -3.01122),
    @(41, 'IMX', 'Immutable', 2.07, 2923933613, 139059990, -3.79748),
    @(42, 'HBAR', 'Hedera', 0.07914400000000001, 2816382211, 186250922, -7.37259),
    @(43, 'KAS', 'Kaspa', 0.118083, 2752302549, 121349364, -6.30125),
    @(44, 'WIF', 'dogwifhat', 2.76, 2724264296, 1113863209, -0.18101),
    @(45, 'MKR', 'Maker', 2919.28, 2692446286, 266140862, -0.533),
    @(46, 'USDE', 'Ethena USDe', 1, 2355595200, 179615642, -0.04662),
    @(47, 'GRT', 'The Graph', 0.241349, 2284072133, 251962864, -9.14095),
    @(48, 'OP', 'Optimism', 2.24, 2237718405, 740545088, -7.01342),
    @(49, 'PEPE', 'Pepe', 0.00000531, 2223933058, 1523553234, -9.520659999999999),
    @(50, 'INJ', 'Injective', 24.74, 2213363943, 478163126, -7.11058),
    @(51, 'XMR', 'Monero', 117.84, 2141065101, 56399940, -5.68921)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
